$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 6125
$ws.Range("I20").Value = 4833.3335
$ws.Range("K20").Value = 4833.3335
$ws.Range("M20").Value = -4603.3335
$ws.Range("H35").Value = 6125
$ws.Range("I35").Value = 4833.3335
$ws.Range("K35").Value = 4833.3335
$ws.Range("M35").Value = -4454.3335
$ws.Range("H138").Value = 3759.7974
$ws.Range("I138").Value = 2420.2354
$ws.Range("J138").Value = 4159.316
$ws.Range("K138").Value = 7260.706200000001
$ws.Range("L138").Value = 12477.948
$ws.Range("M138").Value = -2120.706200000001
$ws.Range("N138").Value = -22757.948
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2318.8333
$ws.Range("I2").Value = 2950
$ws.Range("J2").Value = 2003.25
$ws.Range("K2").Value = 2950
$ws.Range("L2").Value = 2003.25
$ws.Range("M2").Value = -2837
$ws.Range("N2").Value = -2229.25
$ws.Range("H12").Value = 10500
$ws.Range("J12").Value = 10500
$ws.Range("L12").Value = 10500
$ws.Range("N12").Value = -10846
$ws.Range("H32").Value = 6615.202
$ws.Range("I32").Value = 4838.185
$ws.Range("J32").Value = 17687.385
$ws.Range("K32").Value = 4838.185
$ws.Range("L32").Value = 17687.385
$ws.Range("M32").Value = -4551.185
$ws.Range("N32").Value = -18261.385
$ws.Range("H110").Value = 1276.2727
$ws.Range("I110").Value = 892.2857
$ws.Range("J110").Value = 1948.25
$ws.Range("K110").Value = 892.2857
$ws.Range("L110").Value = 1948.25
$ws.Range("M110").Value = 1152.7143
$ws.Range("N110").Value = -6038.25
$ws.Range("H116").Value = 2318.8333
$ws.Range("I116").Value = 2950
$ws.Range("J116").Value = 2003.25
$ws.Range("K116").Value = 2950
$ws.Range("L116").Value = 2003.25
$ws.Range("M116").Value = -656
$ws.Range("N116").Value = -6591.25
$ws.Range("H122").Value = 4035363.8
$ws.Range("I122").Value = 3503.6924
$ws.Range("J122").Value = 6947262.5
$ws.Range("K122").Value = 10511.0772
$ws.Range("L122").Value = 20841787.5
$ws.Range("M122").Value = -8061.0772
$ws.Range("N122").Value = -20846687.5
$ws.Range("H132").Value = 5669.1025
$ws.Range("I132").Value = 1759.3334
$ws.Range("J132").Value = 9020.333000000001
$ws.Range("K132").Value = 5278.0002
$ws.Range("L132").Value = 27060.999
$ws.Range("M132").Value = -2748.0002
$ws.Range("N132").Value = -32120.999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2318.8333
$ws.Range("I3").Value = 2950
$ws.Range("J3").Value = 2003.25
$ws.Range("K3").Value = 2950
$ws.Range("L3").Value = 2003.25
$ws.Range("M3").Value = -2836
$ws.Range("N3").Value = -2231.25
$ws.Range("H107").Value = 2583.3333
$ws.Range("I107").Value = 2583.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2583.3333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -663.3332999999998
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 4152.829
$ws.Range("I134").Value = 3535.75
$ws.Range("J134").Value = 8595.799999999999
$ws.Range("K134").Value = 10607.25
$ws.Range("L134").Value = 25787.4
$ws.Range("M134").Value = -8072.25
$ws.Range("N134").Value = -30857.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1114.1666
$ws.Range("I2").Value = 1199.25
$ws.Range("J2").Value = 944
$ws.Range("K2").Value = 1199.25
$ws.Range("L2").Value = 944
$ws.Range("M2").Value = -1086.25
$ws.Range("N2").Value = -1170
$ws.Range("H31").Value = 1907.0968
$ws.Range("I31").Value = 1437.8955
$ws.Range("J31").Value = 3116.1924
$ws.Range("K31").Value = 1437.8955
$ws.Range("L31").Value = 3116.1924
$ws.Range("M31").Value = -1142.8955
$ws.Range("N31").Value = -3706.1924
$ws.Range("H34").Value = 1907.0968
$ws.Range("I34").Value = 1437.8955
$ws.Range("J34").Value = 3116.1924
$ws.Range("K34").Value = 1437.8955
$ws.Range("L34").Value = 3116.1924
$ws.Range("M34").Value = -1235.8955
$ws.Range("N34").Value = -3520.1924
$ws.Range("H58").Value = 2220556.8
$ws.Range("I58").Value = 3638884.5
$ws.Range("J58").Value = 4419.625
$ws.Range("K58").Value = 3638884.5
$ws.Range("L58").Value = 4419.625
$ws.Range("M58").Value = -3638681.5
$ws.Range("N58").Value = -4825.625
$ws.Range("H107").Value = 1004.8947
$ws.Range("I107").Value = 1008.41174
$ws.Range("J107").Value = 975
$ws.Range("K107").Value = 1008.41174
$ws.Range("L107").Value = 975
$ws.Range("M107").Value = 911.58826
$ws.Range("N107").Value = -4815
$ws.Range("H136").Value = 2220556.8
$ws.Range("I136").Value = 3638884.5
$ws.Range("J136").Value = 4419.625
$ws.Range("K136").Value = 10916653.5
$ws.Range("L136").Value = 13258.875
$ws.Range("M136").Value = -10914103.5
$ws.Range("N136").Value = -18358.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 519.6
$ws.Range("I114").Value = 314.2143
$ws.Range("J114").Value = 998.8333
$ws.Range("K114").Value = 942.6428999999999
$ws.Range("L114").Value = 2996.4999
$ws.Range("M114").Value = 2311.3571
$ws.Range("N114").Value = -9504.499899999999
$ws.Range("H117").Value = 1699
$ws.Range("I117").Value = 352.66666
$ws.Range("J117").Value = 2102.9
$ws.Range("K117").Value = 1057.99998
$ws.Range("L117").Value = 6308.700000000001
$ws.Range("M117").Value = 2384.00002
$ws.Range("N117").Value = -13192.7
$ws.Range("H121").Value = 1371
$ws.Range("J121").Value = 1479.238
$ws.Range("L121").Value = 4437.714
$ws.Range("N121").Value = -7057.714
$ws.Range("H131").Value = 34008.863
$ws.Range("I131").Value = 1642.0834
$ws.Range("J131").Value = 56856
$ws.Range("K131").Value = 4926.2502
$ws.Range("L131").Value = 170568
$ws.Range("M131").Value = 113.7497999999996
$ws.Range("N131").Value = -180648
$ws.Range("H134").Value = 5105.174
$ws.Range("I134").Value = 5598.091
$ws.Range("J134").Value = 4653.3335
$ws.Range("K134").Value = 16794.273
$ws.Range("L134").Value = 13960.0005
$ws.Range("M134").Value = -11724.273
$ws.Range("N134").Value = -24100.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2069.5454
$ws.Range("I113").Value = 2423.3845
$ws.Range("J113").Value = 1558.4445
$ws.Range("K113").Value = 2423.3845
$ws.Range("L113").Value = 1558.4445
$ws.Range("M113").Value = -253.3845000000001
$ws.Range("N113").Value = -5898.4445
$ws.Range("H124").Value = 76450
$ws.Range("J124").Value = 76450
$ws.Range("L124").Value = 76450
$ws.Range("N124").Value = -86270
$ws.Range("H132").Value = 3460.1667
$ws.Range("I132").Value = 3630.1538
$ws.Range("K132").Value = 10890.4614
$ws.Range("M132").Value = -8360.4614
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3549.8333
$ws.Range("I7").Value = 2434.6667
$ws.Range("K7").Value = 2434.6667
$ws.Range("M7").Value = -2322.6667
$ws.Range("H16").Value = 1340
$ws.Range("I16").Value = 1325.4546
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1325.4546
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1155.4546
$ws.Range("N16").Value = -1840
$ws.Range("H40").Value = 4109.4287
$ws.Range("I40").Value = 3849.889
$ws.Range("K40").Value = 3849.889
$ws.Range("M40").Value = -3713.889
$ws.Range("H126").Value = 3549.8333
$ws.Range("I126").Value = 2434.6667
$ws.Range("K126").Value = 7304.000100000001
$ws.Range("M126").Value = -4834.000100000001
$ws.Range("H136").Value = 6705.241
$ws.Range("I136").Value = 4843.467
$ws.Range("J136").Value = 8700
$ws.Range("K136").Value = 14530.401
$ws.Range("L136").Value = 26100
$ws.Range("M136").Value = -11980.401
$ws.Range("N136").Value = -31200
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 41040.25
$ws.Range("J123").Value = 40712.715
$ws.Range("L123").Value = 40712.715
$ws.Range("N123").Value = -50512.715
$ws.Range("H132").Value = 1629.1613
$ws.Range("I132").Value = 1365.5769
$ws.Range("J132").Value = 2999.8
$ws.Range("K132").Value = 4096.7307
$ws.Range("L132").Value = 8999.400000000001
$ws.Range("M132").Value = -1566.7307
$ws.Range("N132").Value = -14059.4

